# Update crypto price/volume table cells per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.036.27"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.953.82"
$ws.Range("E3").Value = "  +0.56%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").Value = "'379.16"
$ws.Range("E5").Value = "  +0.35%  "
$ws.Range("D6").Value = "'102.02"
$ws.Range("E6").Value = "  -0.36%  "
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  -0.11%  "
$ws.Range("D10").Value = "'36.28"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("E11").Value = "  -0.36%  "
$ws.Range("D12").Value = "'0.0849"
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("B13").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C13").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D13").Value = "3.424.90"
$ws.Range("E13").Value = "  +0.60%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'18.43"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("B15").Value = "Uniswap"
$ws.Range("C15").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D15").Value = "'12.42"
$ws.Range("E15").Value = "  +74.45%  "
$ws.Range("D16").Value = "'7.78"
$ws.Range("E16").Value = "  +5.42%  "
$ws.Range("E17").Value = "  +3.80%  "
$ws.Range("D18").Value = "2.955.28"
$ws.Range("E18").Value = "  +0.59%  "
$ws.Range("D19").Value = "51.007.99"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  -3.16%  "
$ws.Range("D21").Value = "'12.37"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("E23").Value = "  +15.92%  "
$ws.Range("D24").Value = "'69.69"
$ws.Range("E24").Value = "  +1.92%  "
$ws.Range("D25").Value = "'266.38"
$ws.Range("E25").Value = "  +1.62%  "
$ws.Range("D26").Value = "'8.08"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("D29").Value = "'25.77"
$ws.Range("E29").Value = "  +0.76%  "
$ws.Range("E30").Value = "  -8.71%  "
$ws.Range("D31").Value = "'0.107"
$ws.Range("E31").Value = "  -6.46%  "
$ws.Range("D32").Value = "'10.33"
$ws.Range("E32").Value = "  +5.43%  "
$ws.Range("D33").Value = "'34.06"
$ws.Range("E33").Value = "  -0.22%  "
$ws.Range("B34").Value = "Toncoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D34").Value = "'2.06"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("B35").Value = "OKB"
$ws.Range("C35").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D35").Value = "'50.60"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "'0.0437"
$ws.Range("E36").Value = "  -3.73%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("E38").Value = "  +6.51%  "
$ws.Range("E39").Value = "  +1.12%  "
$ws.Range("D40").Value = "'16.57"
$ws.Range("E40").Value = "  +0.02%  "
$ws.Range("E41").Value = "  +2.36%  "
$ws.Range("D42").Value = "'2.48"
$ws.Range("E42").Value = "  -3.97%  "
$ws.Range("D43").Value = "'119.93"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "'3.59"
$ws.Range("E44").Value = "  +11.26%  "
$ws.Range("D45").Value = "'21.40"
$ws.Range("E45").Value = "  +0.61%  "
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").Value = "2.017.64"
$ws.Range("E47").Value = "  +0.59%  "
$ws.Range("E48").Value = "  -2.76%  "
$ws.Range("D49").Value = "'0.258"
$ws.Range("E49").Value = "  -6.28%  "
$ws.Range("D50").Value = "'0.0319"
$ws.Range("E50").Value = "  -7.84%  "
$ws.Range("E51").Value = "  +5.48%  "
